$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5294015407562256
$ws.Range("B1").Value = 0.4796834290027618
$ws.Range("C1").Value = 3.511420965194702
$ws.Range("D1").Value = 1.655008554458618
$ws.Range("E1").Value = 1.164329767227173
